# Update gh-pages output (regenerated data scrape) at 456a3b4
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (sheet1): only the "想去人数" (F column) view-counts got
# refreshed for several rows; no rows added/removed.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 91
$ws1.Range("F4").Value = 7244
$ws1.Range("F5").Value = 265
$ws1.Range("F6").Value = 423
$ws1.Range("F7").Value = 3731
$ws1.Range("F9").Value = 534
$ws1.Range("F11").Value = 606
$ws1.Range("F12").Value = 97

# ---------------------------------------------------------------------
# Sheet "演出" (sheet2): the oldest event (Yolo Fes, row 2) has dropped
# off the list; every remaining row shifts up by one, and the "想去人数"
# for the Luke Thompson show ticked up from 4 to 5.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows(2).Delete()
$ws2.Range("F2").Value = 5
# re-number the leading index column after the shift
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2

# ---------------------------------------------------------------------
# Sheet "全部类型" (sheet4): same Yolo Fes row drops off the top, the
# rest shift up by one, and the "想去人数" counters refresh in lock-step
# with the same values used above.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows(2).Delete()
$ws4.Range("F3").Value = 91
$ws4.Range("F4").Value = 5
$ws4.Range("F5").Value = 7244
$ws4.Range("F7").Value = 265
$ws4.Range("F8").Value = 423
$ws4.Range("F9").Value = 3731
$ws4.Range("F11").Value = 534
$ws4.Range("F13").Value = 606
$ws4.Range("F14").Value = 97
# re-number the leading index column after the shift
For ($i = 2; $i -le 14; $i++) {
    $ws4.Range("A$i").Value = $i - 1
}
